# Generate Report for Handback
# Populates the "Latest Target File" (hyperlink) and "Latest Handback File" columns
# for the zh-cn / de-de sheets, refreshes the handback status text + datetimes,
# and widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Colors / formatting used for the existing hyperlink cells (col A) so the new
# hyperlink cells in col I match them visually.
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) == FF6495ED
$hyperlinkUnderline = 2      # xlUnderlineStyleSingle

# Column width helper inputs: the engine stores widths in 1/6-character units
# (raw = round(ColumnWidth*6)/6 + 5/6), so feed it the "characters" width and
# let it snap to the closest supported value.
$wideStatusWidth = 29.1666666666667   # -> ~30 chars (was ~17.2 chars)
$wideFileWidth    = 39.1666666666667  # -> 40 chars (matches col A / col G)

# ---------------------------------------------------------------------------
# Overview sheet: update the status text shown for both locales
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

# Row 2 -> 4d375c79-a9db-48a3-b453-a6097fcdfc18
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", "4d375c79-a9db-48a3-b453-a6097fcdfc18.md") | Out-Null
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("I2").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("J2").Value = "4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 10:53:36"

# Row 3 -> c09d2861-933d-4f7e-a942-dfe71323448b
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", "c09d2861-933d-4f7e-a942-dfe71323448b.md") | Out-Null
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("I3").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("J3").Value = "c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 10:53:36"

$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideFileWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# Row 2 -> 4d375c79-a9db-48a3-b453-a6097fcdfc18
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", "4d375c79-a9db-48a3-b453-a6097fcdfc18.md") | Out-Null
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("I2").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("J2").Value = "4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 10:53:56"

# Row 3 -> c09d2861-933d-4f7e-a942-dfe71323448b
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", "c09d2861-933d-4f7e-a942-dfe71323448b.md") | Out-Null
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("I3").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("J3").Value = "c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 10:53:56"

$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideFileWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth

Write-Host "Handback report generated"
